$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($row, $name, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Update the "last updated" timestamp shown at the top of the report.
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 06:05"

# Refreshed stats for Pakistan (row 22) - country stays in place, numbers updated.
Set-CountryRow 22 "Pakistan" 35788 490 9695 25323 111 9 770

# Rows 174-178: Mongolia moves up in the ranking (new data pushes it above
# Polinesia Francesa / Siria / Angola / Macao, which all shift down one row).
Set-CountryRow 174 "Mongolia" 61 19 15 46 2 0 0
Set-CountryRow 175 "Polinesia Francesa" 60 0 58 2 1 0 0
Set-CountryRow 176 "Siria" 48 0 29 16 0 0 3
Set-CountryRow 177 "Angola" 45 0 14 29 0 0 2
Set-CountryRow 178 "Macao" 45 0 43 2 1 0 0

# Rows 192-194: Nueva Caledonia / Belice move above Santa Lucia.
Set-CountryRow 192 "Nueva Caledonia" 18 0 18 0 0 0 0
Set-CountryRow 193 "Belice" 18 0 16 0 0 0 2
Set-CountryRow 194 "Santa Lucia" 18 0 18 0 0 0 0

# Rows 198-199: Curazao and Dominica swap order.
Set-CountryRow 198 "Curazao" 16 0 14 1 0 0 1
Set-CountryRow 199 "Dominica" 16 0 15 1 0 0 0

# Rows 201-202: Mauritania and Burundi swap order.
Set-CountryRow 201 "Mauritania" 15 0 6 7 0 0 2
Set-CountryRow 202 "Burundi" 15 0 7 7 0 0 1
